$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 6 down to rows 7, 8, and 9 (same data repeated)
$ws.Range("A6:AB6").Copy() | Out-Null
$ws.Range("A7:AB7").PasteSpecial() | Out-Null
$ws.Range("A8:AB8").PasteSpecial() | Out-Null
$ws.Range("A9:AB9").PasteSpecial() | Out-Null
